# Anonymisation des bulletins exemples
# Replace the real student first names on the "Nom" sheet (B4:B19) with
# anonymous letters A-P, and fix the date number format on the "Noel"
# sheet (E1:F1) to dd/mm/yy.

$wb = $excel.ActiveWorkbook

# --- 1. Anonymise student names on the "Nom" sheet ---
$ws = $wb.Worksheets.Item("Nom")
$letters = @("A","B","C","D","E","F","G","H","I","J","K","L","M","N","O","P")
for ($i = 0; $i -lt $letters.Length; $i++) {
    $row = 4 + $i
    $ws.Cells.Item($row, 2).Value = $letters[$i]
}

# --- 2. Fix date format on the "Noel" sheet ---
$wsNoel = $wb.Worksheets.Item("Noel")
$wsNoel.Range("E1:F1").NumberFormat = "dd/mm/yy"
